$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace old email/name with new ones, everywhere they occur
for ($r = 2; $r -le 5; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    if ($cellA.Value2 -eq "anshulgupta1791@gmail.com") {
        $cellA.Value = "anshultest@test.com"
    }
    $cellB = $ws.Cells.Item($r, 2)
    if ($cellB.Value2 -eq "Anshul") {
        $cellB.Value = "AnshulTest"
    }
}

# Update active selection to match post-edit state
$ws.Range("B3").Select()
